$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of C38:C61 (the #DIV/0! error cells caused by the
# shared formula referencing blank cells 12 rows down).
$ws.Range("C38:C61").ClearContents()

# Column B: auto-fit width to content (produces bestFit behaviour).
$ws.Columns("B").AutoFit()
